$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Tabla_Rendimientos")

# Fill the "description" column (C) for all data rows (2-38) with "-"
$rng = $ws.Range("C2:C38")
$rng.Value = "-"
$rng.HorizontalAlignment = -4108  # xlCenter

$rng.Select() | Out-Null
